# The workbook tracks weekly price observations for "Zapallo italiano" at
# Feria Lagunitas de Puerto Montt. A new weekly observation was inserted
# at row 45 (pushing the existing rows 45-151 down to 46-152), dated
# 2021-11-19, with its own volume/price/origin data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45; everything below shifts down one row
# (old row 45 -> new row 46, ..., old row 151 -> new row 152), and the sheet's
# dimension grows from A1:R151 to A1:R152 automatically.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new observation.
$ws.Range("A45").Value = 4
$ws.Range("B45").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C45").Value = "Los Lagos"
$ws.Range("D45").Value = 44519
$ws.Range("E45").Value = 10
$ws.Range("F45").Value = 100112032
$ws.Range("G45").Value = "Zapallo italiano"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 200
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = 12000
$ws.Range("N45").Value = "`$/caja 50 unidades"
$ws.Range("O45").Value = "Región de O'Higgins"
$ws.Range("P45").Value = 240
$ws.Range("Q45").Value = 50
$ws.Range("R45").Value = "Hortaliza"
